$d = $word.ActiveDocument

# Locate the "Requisitos" / "LOQ4056..." paragraph and remove the three
# paragraphs that followed it (an empty paragraph, a page-break-before
# empty paragraph, and the site-footer copyright paragraph), leaving the
# LOQ4056 paragraph directly followed by the remaining empty / page-break
# paragraphs that precede the sectPr.

$marker = "LOQ4056: Química Analítica para Engenharia (Requisito fraco)"

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $marker) {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # Delete the three paragraphs immediately following the marker
    # paragraph (indices target+1, target+1, target+1 since each delete
    # shifts everything up).
    for ($n = 0; $n -lt 3; $n++) {
        $victim = $d.Paragraphs.Item($target + 1)
        $victim.Range.Delete()
    }
}
